$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting (style) from the last existing data row (2784) down through
# the new rows so column A keeps the date/time style (s="1") used throughout
# the table, matching the style of all prior rows.
$ws.Range($ws.Cells.Item(2784,1), $ws.Cells.Item(2784,9)).Copy() | Out-Null
$dest = $ws.Range($ws.Cells.Item(2785,1), $ws.Cells.Item(2815,9))
$dest.PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$excel.CutCopyMode = 0

# Row 2785
$ws.Cells.Item(2785,1).Value2 = 44005.404845162
$ws.Cells.Item(2785,2).Value2 = 'renatafarinon@yahoo.com.br'
$ws.Cells.Item(2785,3).Value2 = 'Hospital Santa Ana - UTI ADULTO - TIPO II'
$ws.Cells.Item(2785,4).Value2 = 10
$ws.Cells.Item(2785,5).Value2 = 0
$ws.Cells.Item(2785,6).Value2 = 3
$ws.Cells.Item(2785,7).Value2 = 0
$ws.Cells.Item(2785,8).Value2 = 0
$ws.Cells.Item(2785,9).Value2 = 'utipoa'

# Row 2786
$ws.Cells.Item(2786,1).Value2 = 44005.4218070139
$ws.Cells.Item(2786,2).Value2 = 'taianivargas@hotmail.com'
$ws.Cells.Item(2786,3).Value2 = 'Hospital Nossa Senhora da Conceição - UTI ADULTO - TIPO III'
$ws.Cells.Item(2786,4).Value2 = 75
$ws.Cells.Item(2786,5).Value2 = 0
$ws.Cells.Item(2786,6).Value2 = 65
$ws.Cells.Item(2786,7).Value2 = 0
$ws.Cells.Item(2786,8).Value2 = 27
$ws.Cells.Item(2786,9).Value2 = 'utipoa'

# Row 2787
$ws.Cells.Item(2787,1).Value2 = 44005.4317244444
$ws.Cells.Item(2787,2).Value2 = 'fernanda.stringhi@maededeus.com.br'
$ws.Cells.Item(2787,3).Value2 = 'Hospital Mãe de Deus - UTI ADULTO - TIPO I'
$ws.Cells.Item(2787,4).Value2 = 60
$ws.Cells.Item(2787,5).Value2 = 0
$ws.Cells.Item(2787,6).Value2 = 52
$ws.Cells.Item(2787,7).Value2 = 9
$ws.Cells.Item(2787,8).Value2 = 4
$ws.Cells.Item(2787,9).Value2 = 'utipoa'

# Row 2788
$ws.Cells.Item(2788,1).Value2 = 44005.4467915625
$ws.Cells.Item(2788,2).Value2 = 'roseuti@gmail.com'
$ws.Cells.Item(2788,3).Value2 = 'Hospital Moinhos de Vento - UTI ADULTO - TIPO III'
$ws.Cells.Item(2788,4).Value2 = 56
$ws.Cells.Item(2788,5).Value2 = 0
$ws.Cells.Item(2788,6).Value2 = 46
$ws.Cells.Item(2788,7).Value2 = 4
$ws.Cells.Item(2788,8).Value2 = 12
$ws.Cells.Item(2788,9).Value2 = 'utipoa'

# Row 2789
$ws.Cells.Item(2789,1).Value2 = 44005.4582543171
$ws.Cells.Item(2789,2).Value2 = 'francojw66@yahoo.com.br'
$ws.Cells.Item(2789,3).Value2 = 'Hospital de Pronto Socorro de Porto Alegre - UTI DE QUEIMADOS'
$ws.Cells.Item(2789,4).Value2 = 4
$ws.Cells.Item(2789,5).Value2 = 0
$ws.Cells.Item(2789,6).Value2 = 4
$ws.Cells.Item(2789,7).Value2 = 0
$ws.Cells.Item(2789,8).Value2 = 0
$ws.Cells.Item(2789,9).Value2 = 'utipoa'

# Row 2790
$ws.Cells.Item(2790,1).Value2 = 44005.4593176736
$ws.Cells.Item(2790,2).Value2 = 'francojw66@yahoo.com.br'
$ws.Cells.Item(2790,3).Value2 = 'Hospital de Pronto Socorro de Porto Alegre - UTI DE QUEIMADOS'
$ws.Cells.Item(2790,4).Value2 = 4
$ws.Cells.Item(2790,5).Value2 = 0
$ws.Cells.Item(2790,6).Value2 = 4
$ws.Cells.Item(2790,7).Value2 = 0
$ws.Cells.Item(2790,8).Value2 = 0
$ws.Cells.Item(2790,9).Value2 = 'utipoa'

# Row 2791
$ws.Cells.Item(2791,1).Value2 = 44005.474418588
$ws.Cells.Item(2791,2).Value2 = 'smarcos@ghc.com.br'
$ws.Cells.Item(2791,3).Value2 = 'Hospital Femina - UTI ADULTO - TIPO II'
$ws.Cells.Item(2791,4).Value2 = 6
$ws.Cells.Item(2791,5).Value2 = 0
$ws.Cells.Item(2791,6).Value2 = 3
$ws.Cells.Item(2791,7).Value2 = 0
$ws.Cells.Item(2791,8).Value2 = 0
$ws.Cells.Item(2791,9).Value2 = 'utipoa'

# Row 2792
$ws.Cells.Item(2792,1).Value2 = 44005.476614919
$ws.Cells.Item(2792,2).Value2 = 'fernanda.ms@santacasa.org.br'
$ws.Cells.Item(2792,3).Value2 = 'Complexo Hospitalar Santa Casa - UTI ADULTO - TIPO III'
$ws.Cells.Item(2792,4).Value2 = 87
$ws.Cells.Item(2792,5).Value2 = 0
$ws.Cells.Item(2792,6).Value2 = 72
$ws.Cells.Item(2792,7).Value2 = 1
$ws.Cells.Item(2792,8).Value2 = 7
$ws.Cells.Item(2792,9).Value2 = 'utipoa'

# Row 2793
$ws.Cells.Item(2793,1).Value2 = 44005.4787971065
$ws.Cells.Item(2793,2).Value2 = 'fernanda.ms@santacasa.org.br'
$ws.Cells.Item(2793,3).Value2 = 'Complexo Hospitalar Santa Casa - UTI PEDIATRICA - TIPO III'
$ws.Cells.Item(2793,4).Value2 = 37
$ws.Cells.Item(2793,5).Value2 = 3
$ws.Cells.Item(2793,6).Value2 = 36
$ws.Cells.Item(2793,7).Value2 = 0
$ws.Cells.Item(2793,8).Value2 = 0
$ws.Cells.Item(2793,9).Value2 = 'utipoa'

# Row 2794
$ws.Cells.Item(2794,1).Value2 = 44005.4838022106
$ws.Cells.Item(2794,2).Value2 = 'marcosboniatti@gmail.com'
$ws.Cells.Item(2794,3).Value2 = 'Hospital Cristo Redentor - UTI ADULTO - TIPO III'
$ws.Cells.Item(2794,4).Value2 = 39
$ws.Cells.Item(2794,5).Value2 = 0
$ws.Cells.Item(2794,6).Value2 = 35
$ws.Cells.Item(2794,7).Value2 = 0
$ws.Cells.Item(2794,8).Value2 = 0
$ws.Cells.Item(2794,9).Value2 = 'utipoa'

# Row 2795
$ws.Cells.Item(2795,1).Value2 = 44005.5014145255
$ws.Cells.Item(2795,2).Value2 = 'renatocvaz@hotmail.com'
$ws.Cells.Item(2795,3).Value2 = 'Instituto de Cardiologia - UTI ADULTO - TIPO III'
$ws.Cells.Item(2795,4).Value2 = 47
$ws.Cells.Item(2795,5).Value2 = 4
$ws.Cells.Item(2795,6).Value2 = 25
$ws.Cells.Item(2795,7).Value2 = 8
$ws.Cells.Item(2795,8).Value2 = 2
$ws.Cells.Item(2795,9).Value2 = 'utipoa'

# Row 2796
$ws.Cells.Item(2796,1).Value2 = 44005.5653323727
$ws.Cells.Item(2796,2).Value2 = 'febueno01@gmail.com'
$ws.Cells.Item(2796,3).Style = 'Normal'
$ws.Cells.Item(2796,4).Value2 = 18
$ws.Cells.Item(2796,5).Value2 = 2
$ws.Cells.Item(2796,6).Value2 = 5
$ws.Cells.Item(2796,7).Value2 = 0
$ws.Cells.Item(2796,8).Value2 = 0
$ws.Cells.Item(2796,9).Value2 = 'utipoa'

# Row 2797
$ws.Cells.Item(2797,1).Value2 = 44005.596984456
$ws.Cells.Item(2797,2).Value2 = 'fredikg@yahoo.com.br'
$ws.Cells.Item(2797,3).Value2 = 'Hospital da Restinga - UTI ADULTO - TIPO II'
$ws.Cells.Item(2797,4).Value2 = 10
$ws.Cells.Item(2797,5).Value2 = 0
$ws.Cells.Item(2797,6).Value2 = 10
$ws.Cells.Item(2797,7).Value2 = 0
$ws.Cells.Item(2797,8).Value2 = 0
$ws.Cells.Item(2797,9).Value2 = 'utipoa'

# Row 2798
$ws.Cells.Item(2798,1).Value2 = 44005.6561149306
$ws.Cells.Item(2798,2).Value2 = 'hlanziotti@gmail.com'
$ws.Cells.Item(2798,3).Value2 = 'Hospital São Lucas - UTI ADULTO - TIPO III'
$ws.Cells.Item(2798,4).Value2 = 59
$ws.Cells.Item(2798,5).Value2 = 0
$ws.Cells.Item(2798,6).Value2 = 46
$ws.Cells.Item(2798,7).Value2 = 3
$ws.Cells.Item(2798,8).Value2 = 8
$ws.Cells.Item(2798,9).Value2 = 'utipoa'

# Row 2799
$ws.Cells.Item(2799,1).Value2 = 44005.6630473032
$ws.Cells.Item(2799,2).Value2 = 'pedrocomerlato@gmail.com'
$ws.Cells.Item(2799,3).Value2 = 'Hospital Independência - UTI ADULTO - TIPO II'
$ws.Cells.Item(2799,4).Value2 = 10
$ws.Cells.Item(2799,5).Value2 = 0
$ws.Cells.Item(2799,6).Value2 = 8
$ws.Cells.Item(2799,7).Value2 = 0
$ws.Cells.Item(2799,8).Value2 = 0
$ws.Cells.Item(2799,9).Value2 = 'utipoa'

# Row 2800
$ws.Cells.Item(2800,1).Value2 = 44005.7331559259
$ws.Cells.Item(2800,2).Value2 = 'ccih@hpa.org'
$ws.Cells.Item(2800,3).Value2 = 'Hospital Porto Alegre - UTI ADULTO - TIPO II'
$ws.Cells.Item(2800,4).Value2 = 7
$ws.Cells.Item(2800,5).Value2 = 0
$ws.Cells.Item(2800,6).Value2 = 3
$ws.Cells.Item(2800,7).Value2 = 0
$ws.Cells.Item(2800,8).Value2 = 0
$ws.Cells.Item(2800,9).Value2 = 'utipoa'

# Row 2801
$ws.Cells.Item(2801,1).Value2 = 44005.7706996528
$ws.Cells.Item(2801,2).Value2 = 'andre.machado@hed.com.br'
$ws.Cells.Item(2801,3).Value2 = 'Hospital Ernesto Dorenelles - UTI ADULTO - TIPO III'
$ws.Cells.Item(2801,4).Value2 = 40
$ws.Cells.Item(2801,5).Value2 = 0
$ws.Cells.Item(2801,6).Value2 = 26
$ws.Cells.Item(2801,7).Value2 = 4
$ws.Cells.Item(2801,8).Value2 = 4
$ws.Cells.Item(2801,9).Value2 = 'utipoa'

# Row 2802
$ws.Cells.Item(2802,1).Value2 = 44005.8009845833
$ws.Cells.Item(2802,2).Value2 = 'smarcos@ghc.com.br'
$ws.Cells.Item(2802,3).Value2 = 'Hospital Femina - UTI ADULTO - TIPO II'
$ws.Cells.Item(2802,4).Value2 = 6
$ws.Cells.Item(2802,5).Value2 = 0
$ws.Cells.Item(2802,6).Value2 = 4
$ws.Cells.Item(2802,7).Value2 = 0
$ws.Cells.Item(2802,8).Value2 = 0
$ws.Cells.Item(2802,9).Value2 = 'utipoa'

# Row 2803
$ws.Cells.Item(2803,1).Value2 = 44005.8623163426
$ws.Cells.Item(2803,2).Value2 = 'joao.krauzer@hmv.org.br'
$ws.Cells.Item(2803,3).Value2 = 'Hospital Moinhos de Vento - UTI PEDIATRICA - TIPO III'
$ws.Cells.Item(2803,4).Value2 = 11
$ws.Cells.Item(2803,5).Value2 = 0
$ws.Cells.Item(2803,6).Value2 = 5
$ws.Cells.Item(2803,7).Value2 = 0
$ws.Cells.Item(2803,8).Value2 = 0
$ws.Cells.Item(2803,9).Value2 = 'utipoa'

# Row 2804
$ws.Cells.Item(2804,1).Value2 = 44006.0229105324
$ws.Cells.Item(2804,2).Value2 = 'taianivargas@hotmail.com'
$ws.Cells.Item(2804,3).Value2 = 'Hospital Nossa Senhora da Conceição - UTI ADULTO - TIPO III'
$ws.Cells.Item(2804,4).Value2 = 75
$ws.Cells.Item(2804,5).Value2 = 0
$ws.Cells.Item(2804,6).Value2 = 65
$ws.Cells.Item(2804,7).Value2 = 0
$ws.Cells.Item(2804,8).Value2 = 27
$ws.Cells.Item(2804,9).Value2 = 'utipoa'

# Row 2805
$ws.Cells.Item(2805,1).Value2 = 44006.3286808449
$ws.Cells.Item(2805,2).Value2 = 'fnagel@hcpa.edu.br'
$ws.Cells.Item(2805,3).Value2 = 'Hospital de Clínicas de Porto Alegre - UTI ADULTO - TIPO III'
$ws.Cells.Item(2805,4).Value2 = 128
$ws.Cells.Item(2805,5).Value2 = 0
$ws.Cells.Item(2805,6).Value2 = 96
$ws.Cells.Item(2805,7).Value2 = 5
$ws.Cells.Item(2805,8).Value2 = 45
$ws.Cells.Item(2805,9).Value2 = 'utipoa'

# Row 2806
$ws.Cells.Item(2806,1).Value2 = 44006.3294286574
$ws.Cells.Item(2806,2).Value2 = 'fnagel@hcpa.edu.br'
$ws.Cells.Item(2806,3).Value2 = 'Hospital de Clínicas de Porto Alegre - UTI ADULTO - TIPO III'
$ws.Cells.Item(2806,4).Value2 = 128
$ws.Cells.Item(2806,5).Value2 = 0
$ws.Cells.Item(2806,6).Value2 = 96
$ws.Cells.Item(2806,7).Value2 = 8
$ws.Cells.Item(2806,8).Value2 = 42
$ws.Cells.Item(2806,9).Value2 = 'utipoa'

# Row 2807
$ws.Cells.Item(2807,1).Value2 = 44006.3535103935
$ws.Cells.Item(2807,2).Value2 = 'cdalmora@hcpa.edu.br'
$ws.Cells.Item(2807,3).Value2 = 'Hospital de Clínicas de Porto Alegre - UTI PEDIATRICA - TIPO III'
$ws.Cells.Item(2807,4).Value2 = 13
$ws.Cells.Item(2807,5).Value2 = 0
$ws.Cells.Item(2807,6).Value2 = 11
$ws.Cells.Item(2807,7).Value2 = 2
$ws.Cells.Item(2807,8).Value2 = 0
$ws.Cells.Item(2807,9).Value2 = 'utipoa'

# Row 2808
$ws.Cells.Item(2808,1).Value2 = 44006.3633799074
$ws.Cells.Item(2808,2).Value2 = 'leandra@portoalegre.rs.gov.br'
$ws.Cells.Item(2808,3).Value2 = 'Hospital Materno Infantil Presidente Vargas - UTI PEDIATRICA - TIPO II'
$ws.Cells.Item(2808,4).Value2 = 12
$ws.Cells.Item(2808,5).Value2 = 0
$ws.Cells.Item(2808,6).Value2 = 3
$ws.Cells.Item(2808,7).Value2 = 0
$ws.Cells.Item(2808,8).Value2 = 0
$ws.Cells.Item(2808,9).Value2 = 'utipoa'

# Row 2809
$ws.Cells.Item(2809,1).Value2 = 44006.3654233681
$ws.Cells.Item(2809,2).Value2 = 'analise.medina@divinaprovidencia.org.br'
$ws.Cells.Item(2809,3).Value2 = 'Hospital Divina Providência - UTI ADULTO - TIPO II'
$ws.Cells.Item(2809,4).Value2 = 20
$ws.Cells.Item(2809,5).Value2 = 0
$ws.Cells.Item(2809,6).Value2 = 19
$ws.Cells.Item(2809,7).Value2 = 2
$ws.Cells.Item(2809,8).Value2 = 1
$ws.Cells.Item(2809,9).Value2 = 'utipoa'

# Row 2810
$ws.Cells.Item(2810,1).Value2 = 44006.3698384143
$ws.Cells.Item(2810,2).Value2 = 'dralubarcellos@gmail.com'
$ws.Cells.Item(2810,3).Value2 = 'Hospital de Pronto Socorro de Porto Alegre - UTI PEDIATRICA - TIPO III'
$ws.Cells.Item(2810,4).Value2 = 8
$ws.Cells.Item(2810,5).Value2 = 0
$ws.Cells.Item(2810,6).Value2 = 5
$ws.Cells.Item(2810,7).Value2 = 0
$ws.Cells.Item(2810,8).Value2 = 0
$ws.Cells.Item(2810,9).Value2 = 'utipoa'

# Row 2811
$ws.Cells.Item(2811,1).Value2 = 44006.3753717593
$ws.Cells.Item(2811,2).Value2 = 'renatafarinon@yahoo.com.br'
$ws.Cells.Item(2811,3).Value2 = 'Hospital Santa Ana - UTI ADULTO - TIPO II'
$ws.Cells.Item(2811,4).Value2 = 10
$ws.Cells.Item(2811,5).Value2 = 0
$ws.Cells.Item(2811,6).Value2 = 5
$ws.Cells.Item(2811,7).Value2 = 0
$ws.Cells.Item(2811,8).Value2 = 0
$ws.Cells.Item(2811,9).Value2 = 'utipoa'

# Row 2812
$ws.Cells.Item(2812,1).Value2 = 44006.3760919792
$ws.Cells.Item(2812,2).Value2 = 'renatafarinon@yahoo.com.br'
$ws.Cells.Item(2812,3).Value2 = 'Hospital Santa Ana - UTI ADULTO - TIPO II'
$ws.Cells.Item(2812,4).Value2 = 10
$ws.Cells.Item(2812,5).Value2 = 0
$ws.Cells.Item(2812,6).Value2 = 5
$ws.Cells.Item(2812,7).Value2 = 0
$ws.Cells.Item(2812,8).Value2 = 0
$ws.Cells.Item(2812,9).Value2 = 'utipoa'

# Row 2813
$ws.Cells.Item(2813,1).Value2 = 44006.3878080208
$ws.Cells.Item(2813,2).Value2 = 'braun.luiz@gmail.com'
$ws.Cells.Item(2813,3).Value2 = 'Hospital Nossa Senhora da Conceição - UTI PEDIATRICA - TIPO II'
$ws.Cells.Item(2813,4).Value2 = 18
$ws.Cells.Item(2813,5).Value2 = 4
$ws.Cells.Item(2813,6).Value2 = 6
$ws.Cells.Item(2813,7).Value2 = 0
$ws.Cells.Item(2813,8).Value2 = 0
$ws.Cells.Item(2813,9).Value2 = 'utipoa'

# Row 2814
$ws.Cells.Item(2814,1).Value2 = 44006.3938132523
$ws.Cells.Item(2814,2).Value2 = 'lufacchi@uol.com.br'
$ws.Cells.Item(2814,3).Value2 = 'Hospital Vila Nova - UTI ADULTO - TIPO II'
$ws.Cells.Item(2814,4).Value2 = 20
$ws.Cells.Item(2814,5).Value2 = 0
$ws.Cells.Item(2814,6).Value2 = 18
$ws.Cells.Item(2814,7).Value2 = 1
$ws.Cells.Item(2814,8).Value2 = 0
$ws.Cells.Item(2814,9).Value2 = 'utipoa'

# Row 2815
$ws.Cells.Item(2815,1).Value2 = 44006.4037261805
$ws.Cells.Item(2815,2).Value2 = 'joao.krauzer@hmv.org.br'
$ws.Cells.Item(2815,3).Value2 = 'Hospital Moinhos de Vento - UTI PEDIATRICA - TIPO III'
$ws.Cells.Item(2815,4).Value2 = 11
$ws.Cells.Item(2815,5).Value2 = 0
$ws.Cells.Item(2815,6).Value2 = 5
$ws.Cells.Item(2815,7).Value2 = 0
$ws.Cells.Item(2815,8).Value2 = 0
$ws.Cells.Item(2815,9).Value2 = 'utipoa'

Write-Output "done"